$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("login")
$wsSignup = $wb.Worksheets.Item("signup")

# ---- login sheet: update the "verify home page header" assertions ----
# The test used to look up the HubSpot / CRMPRO dashboard header; it now
# targets the freeCRM "headertext" element and the signed-in user's name.
$wsLogin.Range("C7").Value = "//td[@class='headertext']"
$wsLogin.Range("E7").Value = "User: Mehraj Ismayilov"
$wsLogin.Range("C8").Value = "headertext"
$wsLogin.Range("E8").Value = "User: Mehraj Ismayilov"

# ---- signup sheet: switch the browser + target URL to freeCRM ----
$wsSignup.Range("E2").Value = "chrome"
$wsSignup.Range("E3").Value = "https://classic.freecrm.com"

# Re-point the hyperlink on E3 at the new freeCRM url (drop the old
# hubspot one first so we don't end up with two hyperlink entries).
$wsSignup.Hyperlinks.Delete()
$wsSignup.Hyperlinks.Add($wsSignup.Range("E3"), "https://classic.freecrm.com")
$wsSignup.Range("E3").Value = "https://classic.freecrm.com"
$wsSignup.Range("E3").Style = "Hyperlink"

# ---- selection / active sheet bookkeeping ----
# Move the lingering selection on the login sheet from E7 to E8, then make
# "signup" the active/visible tab (it was "login" before).
$null = $wsLogin.Range("E8").Select()
$null = $wsSignup.Activate()
$null = $wsSignup.Range("E2").Select()
